# Append two new match rows (137, 138) to the results sheet, matching the
# "Some API changes and UI pushed" commit: two more fixtures were scored by
# the backend and appended to the bottom of the existing results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # The Home xG / Away xG columns store decimal-looking values as TEXT
    # (shared strings) in this workbook, not as numbers. Temporarily mark
    # the cell as Text so Excel doesn't coerce the string into a number,
    # then drop the format back off so no stray cell style lingers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 137: Burnley vs Wolverhampton Wanderers
$ws.Range("A137").Value = "Burnley"
$ws.Range("B137").Value = "Wolverhampton Wanderers"
Set-TextValue $ws.Range("C137") "1.20465"
Set-TextValue $ws.Range("D137") "1.45722"
$ws.Range("E137").Value = 36.4
$ws.Range("F137").Formula = "=SUM(100-E137)"

# Row 138: Chelsea vs West Ham
$ws.Range("A138").Value = "Chelsea"
$ws.Range("B138").Value = "West Ham"
Set-TextValue $ws.Range("C138") "2.41157"
Set-TextValue $ws.Range("D138") "0.277737"
$ws.Range("E138").Value = 53.7
$ws.Range("F138").Formula = "=SUM(100-E138)"

# Match the author's final selection/viewport after the edit.
$ws.Range("F138").Select()
